# "Generate Report for Handback" - mark a.md as handed back and record its
# target/handback file info (Latest Target File / Latest Handback File
# columns) on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Underline + cornflowerblue (FF6495ED in RGB, reversed to BGR for the COM
# Font.Color property) - matches the workbook's existing hyperlink look.
$hyperlinkUnderline = 2
$hyperlinkColor = 15570276

function Style-AsHyperlink($range) {
    $range.Font.Underline = $hyperlinkUnderline
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/aa9fe0a663b7257f7370097613316422c5ab6f42/e2e/a.md"
$bMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/aa9fe0a663b7257f7370097613316422c5ab6f42/e2e/b.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1248db40ecf995a139a6b7a0ec34fd098c504d69/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"

# Status -> handed back
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Handback DateTime
$wsZh.Range("H2").Value = "2016-03-22 20:35:33"
$wsZh.Range("H3").Value = "2016-03-22 20:35:33"

# New "Latest Target File" (F) / "Latest Handback File" (G) cells
$wsZh.Range("F2").Value = "a.md"
$wsZh.Range("G2").Value = $zhXlfName
$wsZh.Range("F3").Value = "a.md"
$wsZh.Range("G3").Value = $zhXlfName

# Rebuild the hyperlinks collection in document order (A2, D2, F2, G2, A3,
# D3, F3, G3) - Hyperlinks.Delete() clears the sheet's whole collection, so
# re-add every link afterwards.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $zhXlfName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $zhXlfName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $bMdUrl, [System.Type]::Missing, [System.Type]::Missing, "b.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $zhXlfName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $zhXlfName) | Out-Null

Style-AsHyperlink $wsZh.Range("F2")
Style-AsHyperlink $wsZh.Range("G2")
Style-AsHyperlink $wsZh.Range("F3")
Style-AsHyperlink $wsZh.Range("G3")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a71031e2ce2360b23c6992329229edfcd33539ce/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfName = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"

# Status -> handed back
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Handback DateTime (distinct timestamp from zh-cn)
$wsDe.Range("H2").Value = "2016-03-22 20:35:44"
$wsDe.Range("H3").Value = "2016-03-22 20:35:44"

# New "Latest Target File" (F) / "Latest Handback File" (G) cells
$wsDe.Range("F2").Value = "a.md"
$wsDe.Range("G2").Value = $deXlfName
$wsDe.Range("F3").Value = "a.md"
$wsDe.Range("G3").Value = $deXlfName

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $deXlfName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $deXlfName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $bMdUrl, [System.Type]::Missing, [System.Type]::Missing, "b.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $deXlfName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $mdUrl, [System.Type]::Missing, [System.Type]::Missing, "a.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, [System.Type]::Missing, [System.Type]::Missing, $deXlfName) | Out-Null

Style-AsHyperlink $wsDe.Range("F2")
Style-AsHyperlink $wsDe.Range("G2")
Style-AsHyperlink $wsDe.Range("F3")
Style-AsHyperlink $wsDe.Range("G3")

Write-Host "Handback report generated."
